$d = $word.ActiveDocument

# 1) Heading text: engagementRessources -> resourcesEngagement
$d.Content.Find.Execute("engagementRessources", $false, $false, $false, $false, $false, $true, 1, $false, "resourcesEngagement", 2) | Out-Null

# 2) "Ressource engagée" -> "Ressource engagée / à engager" (mobilizedResource row, column 2)
$d.Content.Find.Execute("Ressource engagée", $true, $false, $false, $false, $false, $true, 1, $false, "Ressource engagée / à engager", 2) | Out-Null

# 3) Description cell: split into two lines (manual line break = ^l produces <w:br/>)
$oldDescription = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés"
$newDescription = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés en 15-15 et 15-SMUR pour le message RS-RI^lObjet permettant de communiquer la liste des ressources à engager en 15-SMUR pour le message RS-ER"
$d.Content.Find.Execute($oldDescription, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2) | Out-Null

# 4) & 5) Cardinality changes in the "resource" table (2nd table): datetime and vehiculeType rows go from 1..1 to 0..1
$resourceTable = $d.Tables.Item(2)
$resourceTable.Cell(2, 4).Range.Text = "0..1"
$resourceTable.Cell(3, 4).Range.Text = "0..1"
